$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.68826961517334
$ws.Range("B1").Value = 2.073243141174316
$ws.Range("C1").Value = 1.586147546768188
$ws.Range("D1").Value = 1.711413621902466
$ws.Range("E1").Value = 1.51983106136322
